# The two observation records in row 10 ("Garnlav") and row 11
# ("Tretåig hackspett") trade places: every field that differs between
# the rows is swapped, and the public-comment cell (column AC,
# "Ringhack") moves from row 11 to row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writing a plain string like "2026-01-21" into .Value lets Excel's
# input-parsing auto-convert it to a date serial (and stamp a date
# NumberFormat on the cell). Force it to stay literal text, then drop
# the style back to Normal so no stray formatting is left behind.
function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Columns holding plain numbers - safe to swap via .Value2 directly.
$numericCols = @("A", "B", "E", "Q", "R")

# Columns holding text - some look like dates/times, so always route
# them through Set-TextValue to guarantee they stay text.
$textCols = @("F", "G", "H", "Y", "Z", "AA", "AB")

foreach ($col in $numericCols) {
    $addr10 = "$col`10"
    $addr11 = "$col`11"
    $v10 = $ws.Range($addr10).Value2
    $v11 = $ws.Range($addr11).Value2
    $ws.Range($addr10).Value = $v11
    $ws.Range($addr11).Value = $v10
}

foreach ($col in $textCols) {
    $addr10 = "$col`10"
    $addr11 = "$col`11"
    $v10 = $ws.Range($addr10).Value2
    $v11 = $ws.Range($addr11).Value2
    Set-TextValue $ws.Range($addr10) $v11
    Set-TextValue $ws.Range($addr11) $v10
}

# The "Ringhack" public comment (AC) moves from row 11 to row 10;
# row 11's AC cell becomes empty.
Set-TextValue $ws.Range("AC10") "Ringhack"
$ws.Range("AC11").Value = ""
